# Commit: "enabling stochastic increase demand"
#
# On the "Coupling Parameters" sheet:
#  - B31 ("increase demand" switch) is turned ON and becomes a validation
#    formula like its neighbours (mirrors the pattern of C28/C29/...), and
#    picks up the "active switch" highlight formatting that B28/B29 used to
#    have.
#  - B28/B29 lose that highlight (now plain, like B30) since they're no
#    longer the "special" row.
#  - The row 31 label is simplified from "increase demand" to
#    "increase_demand" (matching naming convention of the other keys).
#  - B38 is turned ON.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

# --- capture/move the "active switch" highlight formatting -----------------
# B28 currently carries the highlighted style; grab it before B28 is
# reformatted below, and stamp it onto B31 (the cell becoming the active
# switch).
$ws.Range("B28").Copy()
[void]$ws.Range("B31").PasteSpecial(-4122)   # xlPasteFormats

# B28 and B29 drop the highlight, taking on the plain style used by B30.
$ws.Range("B30").Copy()
[void]$ws.Range("B28").PasteSpecial(-4122)   # xlPasteFormats
[void]$ws.Range("B29").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = $false

# --- row 31: label, value, and formula --------------------------------------
$ws.Range("A31").Value = "increase_demand"
$ws.Range("B31").Value = $true
$ws.Range("C31").Formula = '=IF(OR(B28<>TRUE,B29<>TRUE),"demand and profiles must be fix!!!!!!!!!!","ok")'

# --- row 38: enable the flag -------------------------------------------------
$ws.Range("B38").Value = $true

# --- restore view/selection state -------------------------------------------
[void]$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("B31").Select()

$wb.Application.Calculate()
